$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text so formatted numbers (e.g. "1.00", "0.999")
# are preserved literally instead of being normalized as numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.298.05"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.382.01"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.05"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.19"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.382.31"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.959.65"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.99"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.411.81"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.418.33"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.03"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.54"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.18"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.529.14"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.61"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.412.31"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.89"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "167.69"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.91"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.81"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.779"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.521.59"

# Volume(1h) (column E) updates - percentage strings with surrounding padding spaces.
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("E6").Value = "  -6.44%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -6.99%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("E33").Value = "  -7.96%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("E40").Value = "  -6.97%  "
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -3.82%  "
